# Applies the "documents updated for current day" edit to the Team 04
# project-meeting minutes document.
#
# Summary of changes (see commit message / diff):
#  1. Purpose: "Revision of requirements" -> "Mid-term presentation" (two runs)
#  2. Day/Date/Time: "Tuesday, June 16" -> "Wednesday, June 17"
#  3. Attendee names: re-split into runs w/ spell-check proofErr markers
#  4. Agenda bullet 1: "Revision of requirements" -> "Discuss Mid-term presentation"
#  5. Agenda bullet 2: "Working on clients request mentioned in meeting"
#     -> "Presentation preparation" (three runs) + _GoBack bookmark
#  6. Action item 1 text -> "Work on slides for mid-term presentation"
#  7. Action item 1 due date: 06/16 -> 06/17
#  8. Action item row 2 ("Gather revised requirements from client") removed
#  9. Next meeting date: "Wednesday" + bookmark + ", June 17" -> "Thursday, June 18"
#     (bookmark relocates to change #5 above, matching where Word would
#     leave the _GoBack mark after the most recent edit)

$d = $word.ActiveDocument

$wdXmlPart = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">{0}</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($paragraph, [string]$innerBodyXml) {
    $full = [string]::Format($wdXmlPart, "<w:body>$innerBodyXml</w:body>")
    $paragraph.Range.InsertXML($full)
}

# ---------------------------------------------------------------------
# 1. Purpose cell: "Revision of requirements" -> "Mid-" + "term presentation"
# ---------------------------------------------------------------------
$tbl1 = $d.Tables.Item(1)
$purposePara = $tbl1.Cell(2, 2).Range.Paragraphs.Item(1)
Set-ParagraphXml $purposePara @'
<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Mid-</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>term presentation</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 2. Day, Date & Time cell: "Tuesday, June 16" -> "Wednesday, June 17"
#    (keep the ", 2020 4:00 PM CST" run separate, as in the target XML)
# ---------------------------------------------------------------------
$dateTimePara = $tbl1.Cell(3, 2).Range.Paragraphs.Item(1)
Set-ParagraphXml $dateTimePara @'
<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Wednesday, June 17</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>, 2020 4:00 PM CST</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 3. Attendees: re-split each name into proofed runs (text unchanged)
# ---------------------------------------------------------------------
$attendeesCell = $tbl1.Cell(5, 2)

$attendeePara1 = $attendeesCell.Range.Paragraphs.Item(2)
Set-ParagraphXml $attendeePara1 @'
<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Bhavya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Deepthi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Gorrepati</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

$attendeePara2 = $attendeesCell.Range.Paragraphs.Item(3)
Set-ParagraphXml $attendeePara2 @'
<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Mahalakshmi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Kongari</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

$attendeePara3 = $attendeesCell.Range.Paragraphs.Item(4)
Set-ParagraphXml $attendeePara3 @'
<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Sai Jyothsna </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Mathi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

$attendeePara4 = $attendeesCell.Range.Paragraphs.Item(5)
Set-ParagraphXml $attendeePara4 @'
<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Jeevan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Reddy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Mure</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

$attendeePara5 = $attendeesCell.Range.Paragraphs.Item(6)
Set-ParagraphXml $attendeePara5 @'
<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Dheeraj</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Edupuganti</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

# ---------------------------------------------------------------------
# 4. Agenda bullet 1: "Revision of requirements" -> "Discuss Mid-term presentation"
# ---------------------------------------------------------------------
$agenda1Rng = $d.Content
$agenda1Rng.Find.ClearFormatting()
$found = $agenda1Rng.Find.Execute("Revision of requirements", $true, $false, $false, $false, $false, $true, 1, $false, "Discuss Mid-term presentation", 1)

# ---------------------------------------------------------------------
# 5. Agenda bullet 2: "Working on clients request mentioned in meeting"
#    -> "Presentation " + "prepa" + "ration" runs, plus the _GoBack bookmark
# ---------------------------------------------------------------------
$agenda2Para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Working on clients request*") {
        $agenda2Para = $cand
        break
    }
}
Set-ParagraphXml $agenda2Para @'
<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Presentation </w:t></w:r><w:r><w:t>prepa</w:t></w:r><w:r><w:t>ration</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

# ---------------------------------------------------------------------
# 6 & 7. Action item table: update task text and due-date day
# ---------------------------------------------------------------------
$tbl2 = $d.Tables.Item(2)
$taskPara = $tbl2.Cell(2, 1).Range.Paragraphs.Item(1)
Set-ParagraphXml $taskPara @'
<w:p><w:r><w:t>Work on slides for mid-term presentation</w:t></w:r></w:p>
'@

$dueDatePara = $tbl2.Cell(2, 3).Range.Paragraphs.Item(1)
Set-ParagraphXml $dueDatePara @'
<w:p><w:r><w:t>06/17</w:t></w:r><w:r><w:t>/2020</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 8. Remove the second action-item row entirely
# ---------------------------------------------------------------------
$tbl2.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# 9. Next meeting date: "Wednesday" + bookmark + ", June 17" -> "Thursday, June 18"
# ---------------------------------------------------------------------
$nextMeetingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Date and time of next project team meeting*") {
        $nextMeetingPara = $cand
        break
    }
}
Set-ParagraphXml $nextMeetingPara @'
<w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Date and time of next project team meeting: </w:t></w:r><w:r><w:t>Thursday, June 18</w:t></w:r><w:r><w:t>, 2020 4:00 PM CST</w:t></w:r></w:p>
'@

Write-Output "All edits applied."
